# "change excel2json to many" — update the risk-monitor workbook:
#  - sheet 總覽 (Overview): P/C ratio + margin totals refreshed
#  - sheet 詳細數據 (Detail): mirrored totals refreshed
#  - sheet 個股籌碼 (Per-stock chips): the per-broker detail columns
#    (P..U) are dropped for every stock row and the data-source flag (V)
#    flips from "partial" to "N/A"; a few stocks also got revised
#    price/volume figures (C/D/E).

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 總覽 (Overview) ----------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C11").Value = "'158.39%"
$ws1.Range("E11").Value = "'171.77%"
$ws1.Range("G14").Value = "-26.2億"
$ws1.Range("H14").Value = "-523.98億"

# ---- Sheet 2: 詳細數據 (Detail) --------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B16").Value = "-26.2億"
$ws2.Range("B17").Value = "-523.98億"
$ws2.Range("B21").Value = "'171.77%"

# ---- Sheet 3: 個股籌碼 (Per-stock chips) ------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Revised price/volume figures for a handful of stocks.
$ws3.Range("C15").Value = 1000
$ws3.Range("D15").Value = -2.44
$ws3.Range("E15").Value = 905

$ws3.Range("C16").Value = 306.5
$ws3.Range("D16").Value = -8.1
$ws3.Range("D16").Font.Color = $ws3.Range("D17").Font.Color
$ws3.Range("E16").Value = 26486

$ws3.Range("C17").Value = 132.5
$ws3.Range("D17").Value = -2.57
$ws3.Range("E17").Value = 1591

$ws3.Range("C18").Value = 318.5
$ws3.Range("D18").Value = -0.16
$ws3.Range("E18").Value = 20497

# Every stock row (4..19) loses its per-broker detail columns (P..U) and
# its data-source marker flips from "partial" to "N/A".
for ($row = 4; $row -le 19; $row++) {
    $ws3.Range("P$row`:U$row").Clear()
    $ws3.Range("V$row").Value = "N/A"
}
